$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A178").Value = 177
$ws.Range("B178").Value = 1
$ws.Range("C178").Value = "2024-06-18 17:12:27"
$ws.Range("D178").Value = 200
$ws.Range("E178").Value = 17

$ws.Range("A179").Value = 178
$ws.Range("B179").Value = 2
$ws.Range("C179").Value = "2024-06-18 17:12:27"
$ws.Range("D179").Value = 200
$ws.Range("E179").Value = 3
